$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LEGALLABELLING")

# 1. "CE mark" -> "UKCA mark" (row 10, column D)
$ws.Range("D10").Value = "the information includes the UKCA mark"

# 2. Remove the "the information includes the CE representative details" row
#    (row 12: And / the information includes the CE representative details)
$ws.Rows("12:12").Delete()
